$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 877, shifting the existing row 877 (and all
# rows below it) down by one. This matches the diff where a new data row
# (2026/02/26, 木, 9, 201) is inserted just before the old "2026/12/29"
# block, and the sheet's used range grows from A1:D918 to A1:D919.
$ws.Rows.Item(877).Insert()

# Column A holds the date as a literal text string (e.g. "2026/02/26"),
# not a real date value. Setting a date-looking string via .Value would
# normally get auto-converted into a date serial by Excel, so temporarily
# mark the cell as Text before assigning it, then restore its style back
# to the workbook default so no stray formatting is left behind.
$ws.Range("A877").NumberFormat = "@"
$ws.Range("A877").Value = "2026/02/26"
$ws.Range("A877").Style = "Normal"

$ws.Range("B877").Value = "木"
$ws.Range("C877").Value = 9
$ws.Range("D877").Value = 201
